# The sheet had a stray, data-less row ("grandes regiões e unidades da
# federação") inserted at row 6, which pushed every region's values one row
# below its label (e.g. "norte" in A7 showed the values that actually belong
# to "rondônia", etc.). Deleting that row fixes the data/label alignment:
# rows 7..37 shift up to 6..36, and the now-unused shared string is dropped
# automatically when the workbook is saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("6").Delete()
